$wb = $excel.ActiveWorkbook

# --- Grocery sheet ---
$ws = $wb.Worksheets.Item("Grocery")
$ws.Range("C1").Value = "Price"
$ws.Range("C2").Value = 11.99
$ws.Range("C3").Value = 7
$ws.Range("C4").Value = 5
$ws.Range("C5").Value = 9.99
$ws.Range("C6").Value = 8.99
$ws.Range("C7").Value = 8.99
$ws.Range("C8").Value = 6.5

# --- Pet sheet ---
$ws = $wb.Worksheets.Item("Pet")
$ws.Range("C1").Value = "Price"
$ws.Range("C2").Value = 4.99
$ws.Range("C3").Value = 11.99
$ws.Range("C4").Value = 11.99
$ws.Range("C5").Value = 9.99

# --- Bath sheet ---
$ws = $wb.Worksheets.Item("Bath")
$ws.Range("C1").Value = "Price"
$ws.Range("C2").Value = 15
$ws.Range("C3").Value = 7
$ws.Range("C4").Value = 9
$ws.Range("C5").Value = 6
$ws.Range("C6").Value = 12.5

# Selections matching final state (column C, one row below data) on each sheet
$wb.Worksheets.Item("Grocery").Range("C9").Select() | Out-Null
$wb.Worksheets.Item("Pet").Range("C5").Select() | Out-Null

# Bath becomes the active/selected sheet, with selection at C7
$wsBath = $wb.Worksheets.Item("Bath")
$wsBath.Activate()
$wsBath.Range("C7").Select() | Out-Null
